$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1511.875
$ws.Range("I6").Value = 1511.875
$ws.Range("K6").Value = 4535.625
$ws.Range("M6").Value = -4423.625

$ws.Range("H8").Value = 131.16667
$ws.Range("I8").Value = 137.4
$ws.Range("J8").Value = 100
$ws.Range("K8").Value = 412.2
$ws.Range("L8").Value = 300
$ws.Range("M8").Value = -273.2
$ws.Range("N8").Value = -578

$ws.Range("H33").Value = 66831.60000000001
$ws.Range("I33").Value = 71533.86
$ws.Range("J33").Value = 1000
$ws.Range("K33").Value = 71533.86
$ws.Range("L33").Value = 1000
$ws.Range("M33").Value = -71304.86
$ws.Range("N33").Value = -1458

$ws.Range("H34").Value = 18584
$ws.Range("I34").Value = 17617.6
$ws.Range("J34").Value = 21000
$ws.Range("K34").Value = 17617.6
$ws.Range("L34").Value = 21000
$ws.Range("M34").Value = -17414.6
$ws.Range("N34").Value = -21406

$ws.Range("H36").Value = 18584
$ws.Range("I36").Value = 17617.6
$ws.Range("J36").Value = 21000
$ws.Range("K36").Value = 17617.6
$ws.Range("L36").Value = 21000
$ws.Range("M36").Value = -16902.6
$ws.Range("N36").Value = -22430

$ws.Range("H107").Value = 588666.6
$ws.Range("I107").Value = 714643.1
$ws.Range("J107").Value = 776.3333
$ws.Range("K107").Value = 714643.1
$ws.Range("L107").Value = 776.3333
$ws.Range("M107").Value = -712723.1
$ws.Range("N107").Value = -4616.3333

$ws.Range("H121").Value = 836.913
$ws.Range("J121").Value = 840.4286
$ws.Range("L121").Value = 2521.2858
$ws.Range("N121").Value = -6015.2858

$ws.Range("H132").Value = 1900
$ws.Range("I132").Value = 1517.0731
$ws.Range("J132").Value = 3327.2727
$ws.Range("K132").Value = 4551.219300000001
$ws.Range("L132").Value = 9981.8181
$ws.Range("M132").Value = -2021.219300000001
$ws.Range("N132").Value = -15041.8181

$ws.Range("H137").Value = 7577597.5
$ws.Range("I137").Value = 1504.1951
$ws.Range("J137").Value = 20002390
$ws.Range("K137").Value = 4512.5853
$ws.Range("L137").Value = 60007170
$ws.Range("M137").Value = -1962.5853
$ws.Range("N137").Value = -60012270

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19315.598
$ws.Range("I32").Value = 16740.688
$ws.Range("K32").Value = 16740.688
$ws.Range("M32").Value = -16453.688

$ws.Range("H61").Value = 1365.8975
$ws.Range("I61").Value = 1199.5
$ws.Range("K61").Value = 1199.5
$ws.Range("M61").Value = -987.5

$ws.Range("H132").Value = 2325.5
$ws.Range("I132").Value = 1871.4
$ws.Range("J132").Value = 2893.125
$ws.Range("K132").Value = 5614.200000000001
$ws.Range("L132").Value = 8679.375
$ws.Range("M132").Value = -3084.200000000001
$ws.Range("N132").Value = -13739.375

$ws.Range("H136").Value = 1365.8975
$ws.Range("I136").Value = 1199.5
$ws.Range("K136").Value = 3598.5
$ws.Range("M136").Value = -1048.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 45487264
$ws.Range("I20").Value = 75344.125
$ws.Range("K20").Value = 75344.125
$ws.Range("M20").Value = -75097.125

$ws.Range("H134").Value = 59462.484
$ws.Range("I134").Value = 2499.5186
$ws.Range("J134").Value = 251712.5
$ws.Range("K134").Value = 7498.5558
$ws.Range("L134").Value = 755137.5
$ws.Range("M134").Value = -4963.5558
$ws.Range("N134").Value = -760207.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 666.6667
$ws.Range("I8").Value = 500
$ws.Range("J8").Value = 750
$ws.Range("K8").Value = 500
$ws.Range("L8").Value = 750
$ws.Range("M8").Value = -360
$ws.Range("N8").Value = -1030

$ws.Range("H31").Value = 1903.4565
$ws.Range("I31").Value = 1092.0834
$ws.Range("K31").Value = 1092.0834
$ws.Range("M31").Value = -797.0834

$ws.Range("H34").Value = 1903.4565
$ws.Range("I34").Value = 1092.0834
$ws.Range("K34").Value = 1092.0834
$ws.Range("M34").Value = -890.0834

$ws.Range("H132").Value = 2325.4
$ws.Range("I132").Value = 750.8570999999999
$ws.Range("J132").Value = 5999.3335
$ws.Range("K132").Value = 2252.5713
$ws.Range("L132").Value = 17998.0005
$ws.Range("M132").Value = 277.4287000000004
$ws.Range("N132").Value = -23058.0005

$ws.Range("H134").Value = 2707.2974
$ws.Range("I134").Value = 2106.0454
$ws.Range("J134").Value = 3589.1333
$ws.Range("K134").Value = 6318.1362
$ws.Range("L134").Value = 10767.3999
$ws.Range("M134").Value = -3783.1362
$ws.Range("N134").Value = -15837.3999

$ws.Range("H140").Value = 50953.332
$ws.Range("J140").Value = 50953.332
$ws.Range("L140").Value = 50953.332
$ws.Range("N140").Value = -61313.332

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 5900.8335
$ws.Range("J51").Value = 5900.8335
$ws.Range("L51").Value = 17702.5005
$ws.Range("N51").Value = -18622.5005

$ws.Range("H68").Value = 1041.1235
$ws.Range("I68").Value = 585.95654
$ws.Range("J68").Value = 1528.0465
$ws.Range("K68").Value = 1757.86962
$ws.Range("L68").Value = 4584.139499999999
$ws.Range("M68").Value = -946.8696199999999
$ws.Range("N68").Value = -6206.139499999999

$ws.Range("H71").Value = 1041.1235
$ws.Range("I71").Value = 585.95654
$ws.Range("J71").Value = 1528.0465
$ws.Range("K71").Value = 5273.60886
$ws.Range("L71").Value = 13752.4185
$ws.Range("M71").Value = -1217.60886
$ws.Range("N71").Value = -21864.4185

$ws.Range("H92").Value = 30431.9
$ws.Range("I92").Value = 60225.2
$ws.Range("J92").Value = 638.6
$ws.Range("K92").Value = 180675.6
$ws.Range("L92").Value = 1915.8
$ws.Range("M92").Value = -179427.6
$ws.Range("N92").Value = -4411.8

$ws.Range("H122").Value = 400462.47
$ws.Range("I122").Value = 410.1579
$ws.Range("J122").Value = 1667294.9
$ws.Range("K122").Value = 3691.4211
$ws.Range("L122").Value = 15005654.1
$ws.Range("M122").Value = -1241.4211
$ws.Range("N122").Value = -15010554.1

$ws.Range("H124").Value = 3327.7778
$ws.Range("I124").Value = 1237.5
$ws.Range("K124").Value = 3712.5
$ws.Range("M124").Value = 1197.5

$ws.Range("H125").Value = 2749.8333
$ws.Range("I125").Value = 1500
$ws.Range("J125").Value = 2999.8
$ws.Range("K125").Value = 4500
$ws.Range("L125").Value = 8999.400000000001
$ws.Range("M125").Value = 420
$ws.Range("N125").Value = -18839.4

$ws.Range("H136").Value = 3432
$ws.Range("I136").Value = 2560.889
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 7682.667
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -2582.667
$ws.Range("N136").Value = -25200

$ws.Range("H138").Value = 2727.6924
$ws.Range("I138").Value = 2147
$ws.Range("K138").Value = 6441
$ws.Range("M138").Value = -1301

$ws.Range("H139").Value = 1794.9286
$ws.Range("I139").Value = 1384.4546
$ws.Range("J139").Value = 3300
$ws.Range("K139").Value = 4153.3638
$ws.Range("L139").Value = 9900
$ws.Range("M139").Value = 986.6361999999999
$ws.Range("N139").Value = -20180

$ws.Range("H140").Value = 4661
$ws.Range("I140").Value = 5026.857
$ws.Range("K140").Value = 15080.571
$ws.Range("M140").Value = -9900.571

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 9745
$ws.Range("J33").Value = 9745
$ws.Range("L33").Value = 9745
$ws.Range("N33").Value = -10249

$ws.Range("H80").Value = 4333.8335
$ws.Range("I80").Value = 4722.222
$ws.Range("K80").Value = 4722.222
$ws.Range("M80").Value = -3724.222

$ws.Range("H83").Value = 4333.8335
$ws.Range("I83").Value = 4722.222
$ws.Range("K83").Value = 23611.11
$ws.Range("M83").Value = -18619.11

$ws.Range("H138").Value = 39738.89
$ws.Range("J138").Value = 39738.89
$ws.Range("L138").Value = 39738.89
$ws.Range("N138").Value = -50018.89

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 466.66666
$ws.Range("I46").Value = 425
$ws.Range("J46").Value = 487.5
$ws.Range("K46").Value = 425
$ws.Range("L46").Value = 487.5
$ws.Range("M46").Value = -237
$ws.Range("N46").Value = -863.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1566.8462
$ws.Range("I81").Value = 1579
$ws.Range("J81").Value = 1500
$ws.Range("K81").Value = 3158
$ws.Range("L81").Value = 3000
$ws.Range("M81").Value = -2097
$ws.Range("N81").Value = -5122

$ws.Range("H84").Value = 1566.8462
$ws.Range("I84").Value = 1579
$ws.Range("J84").Value = 1500
$ws.Range("K84").Value = 15790
$ws.Range("L84").Value = 15000
$ws.Range("M84").Value = -10486
$ws.Range("N84").Value = -25608

$ws.Range("H136").Value = 2062.7693
$ws.Range("I136").Value = 2310.5405
$ws.Range("K136").Value = 6931.6215
$ws.Range("M136").Value = -4381.6215
